$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4, pushing existing rows 4..34 down to 5..35.
$ws.Rows(4).Insert()

# Populate the new row 4 with the new weekly record (same market/category
# metadata as the surrounding rows, new date + volume + price figures).
$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Vega Monumental Concepción"
$ws.Range("C4").Value = "Bíobío"
$ws.Range("D4").Value = "2022-07-19"
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 100114007
$ws.Range("G4").Value = "Jengibre"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 14400
$ws.Range("N4").Value = "$/caja 13 kilos"
$ws.Range("O4").Value = "Perú"
$ws.Range("P4").Value = 1108
$ws.Range("Q4").Value = 13
$ws.Range("R4").Value = "Hortaliza"
